$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) for rows 2-14 changed from serial date 45190 to 45192
for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 3).Value2 = 45192
}
